# Weekly update: a new price record is published for the most recent week,
# pushing every existing record down by one row (row 11 -> 12, ... row 73 -> 74)
# and inserting the new record at row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 11-73 down to 12-74 by inserting a new row at 11.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the latest week's record.
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 44819
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112009
$ws.Range("G11").Value = "Acelga"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 1800
$ws.Range("M11").Value = 1650
$ws.Range("N11").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 550
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = "Hortaliza"
